$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from column K (rows 4-14) into the new column L
$ws.Range("K4:K14").Copy()
$ws.Range("L4:L14").PasteSpecial(-4122)

# Populate the new 2023 data column (L)
$ws.Range("L4").Value = 2023
$ws.Range("L5").Value = 22.743990309495757
$ws.Range("L6").Value = 52.401334422687093
$ws.Range("L7").Value = 40.084286291781751
$ws.Range("L8").Value = 58.6564425462321
$ws.Range("L9").Value = 52.689880705632987
$ws.Range("L10").Value = 19.88866894869804
$ws.Range("L11").Value = 35.972443863264772
$ws.Range("L12").Value = 12.061786277026036
$ws.Range("L13").Value = -0.064288010286095529
$ws.Range("L14").Value = 34.132731805770057

# Adjust row heights to match the updated layout
$ws.Rows.Item(1).RowHeight = 67.5
$ws.Rows.Item(4).RowHeight = 14.25
$ws.Rows.Item(5).RowHeight = 14.25
$ws.Rows.Item(6).RowHeight = 14.25
$ws.Rows.Item(7).RowHeight = 14.25
$ws.Rows.Item(8).RowHeight = 14.25
$ws.Rows.Item(9).RowHeight = 14.25
$ws.Rows.Item(10).RowHeight = 14.25
$ws.Rows.Item(11).RowHeight = 14.25
$ws.Rows.Item(12).RowHeight = 14.25
$ws.Rows.Item(13).RowHeight = 14.25
$ws.Rows.Item(14).RowHeight = 14.25

# Reset the saved selection back to the default (A1), clearing the stray
# "M7" selection that was stored in the sheet view
$ws.Range("A1").Select()
